# Apply updated K (strikeout) values to the save_data sheet for pruitt_austin.
# This regenerates column G ("K", formerly "Strike#") with corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 0
    4  = 1
    5  = 2
    7  = 1
    8  = 2
    9  = 0
    10 = 0
    11 = 1
    12 = 0
    13 = 3
    14 = 2
    15 = 0
    16 = 0
    17 = 1
    18 = 2
    19 = 1
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 1
    26 = 1
    27 = 1
    28 = 1
    29 = 1
    30 = 0
    31 = 0
    32 = 1
    33 = 2
    34 = 0
    35 = 1
    36 = 2
    37 = 2
    38 = 0
    39 = 0
    41 = 2
    42 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
